# Corrections following third round of review
#
# Remove the "subgenus" column from the Materials sheet (header
# "subgenus" in row 1, value "${subgenus}" in row 2). Deleting the
# entire column shifts every subsequent column one position to the
# left and lets the workbook drop the now-unused "subgenus" /
# "${subgenus}" shared-string entries automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

$ws.Range("AS1:AS2").EntireColumn.Delete()
